$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "22.065.67"
$ws.Range("E2").Value = "  -0.78%  "
$ws.Range("D3").Value = "1.555.98"
$ws.Range("E3").Value = "  -0.17%  "
$ws.Range("D4").Value = "'0.9995"
$ws.Range("E4").Value = "  -0.12%  "
$ws.Range("E5").Value = "  -0.03%  "
$ws.Range("D6").Value = "'290.67"
$ws.Range("E6").Value = "  +0.79%  "
$ws.Range("D7").Value = "'0.3955"
$ws.Range("E7").Value = "  +3.73%  "
$ws.Range("D8").Value = "'0.3225"
$ws.Range("E8").Value = "  -2.95%  "
$ws.Range("D9").Value = "'43.83"
$ws.Range("E9").Value = "  -2.09%  "
$ws.Range("D10").Value = "'0.07266"
$ws.Range("E10").Value = "  -1.95%  "
$ws.Range("D11").Value = "'1.079"
$ws.Range("E11").Value = "  -5.73%  "
$ws.Range("D12").Value = "'0.9996"
$ws.Range("E12").Value = "  -0.13%  "
$ws.Range("D13").Value = "'5.701"
$ws.Range("E13").Value = "  -2.56%  "
$ws.Range("D14").Value = "'18.83"
$ws.Range("E14").Value = "  -6.95%  "
$ws.Range("D15").Value = "'0.00001133"
$ws.Range("E15").Value = "  +5.28%  "
$ws.Range("D16").Value = "'6.639"
$ws.Range("E16").Value = "  -1.70%  "
$ws.Range("D17").Value = "1.554.85"
$ws.Range("E17").Value = "  -0.69%  "
$ws.Range("D18").Value = "'0.06595"
$ws.Range("E18").Value = "  -0.99%  "
$ws.Range("D19").Value = "'83.62"
$ws.Range("E19").Value = "  -3.33%  "
$ws.Range("D20").Value = "'1.000"
$ws.Range("E20").Value = "  -0.03%  "
$ws.Range("D21").Value = "'6.288"
$ws.Range("E21").Value = "  -1.95%  "
$ws.Range("D22").Value = "'15.56"
$ws.Range("E22").Value = "  -3.79%  "
$ws.Range("D23").Value = "'11.33"
$ws.Range("E23").Value = "  -3.58%  "
$ws.Range("D24").Value = "22.073.51"
$ws.Range("E24").Value = "  -0.68%  "
$ws.Range("E25").Value = "  +3.89%  "
$ws.Range("D26").Value = "'2.429"
$ws.Range("E26").Value = "  -5.13%  "
$ws.Range("D27").Value = "'148.87"
$ws.Range("E27").Value = "  -1.60%  "
$ws.Range("D28").Value = "'18.69"
$ws.Range("E28").Value = "  -3.16%  "
$ws.Range("D29").Value = "'4.884"
$ws.Range("E29").Value = "  -1.10%  "
$ws.Range("D30").Value = "1.730.11"
$ws.Range("E30").Value = "  -0.35%  "
$ws.Range("D31").Value = "'119.15"
$ws.Range("E31").Value = "  -3.34%  "
$ws.Range("D32").Value = "'0.9773"
$ws.Range("E32").Value = "  -10.58%  "
$ws.Range("D33").Value = "'5.849"
$ws.Range("E33").Value = "  -1.17%  "
$ws.Range("D34").Value = "'0.08331"
$ws.Range("E34").Value = "  +1.28%  "
$ws.Range("D35").Value = "'9.127"
$ws.Range("E35").Value = "  -2.14%  "
$ws.Range("D36").Value = "'1.603"
$ws.Range("E36").Value = "  -16.01%  "
$ws.Range("D37").Value = "'0.02270"
$ws.Range("E37").Value = "  -2.72%  "
$ws.Range("D38").Value = "'5.127"
$ws.Range("E38").Value = "  -3.68%  "
$ws.Range("D39").Value = "'0.06000"
$ws.Range("E39").Value = "  -5.23%  "
$ws.Range("D40").Value = "'1.209"
$ws.Range("E40").Value = "  -1.91%  "
$ws.Range("D41").Value = "'0.2038"
$ws.Range("E41").Value = "  -5.88%  "
$ws.Range("E42").Value = "  +0.00%  "
$ws.Range("D43").Value = "'10.72"
$ws.Range("E43").Value = "  -2.64%  "
$ws.Range("D44").Value = "'0.5829"
$ws.Range("E44").Value = "  -4.01%  "
$ws.Range("D45").Value = "'3.747"
$ws.Range("E45").Value = "  +0.04%  "
$ws.Range("D46").Value = "'12.95"
$ws.Range("E46").Value = "  -5.99%  "
$ws.Range("D47").Value = "'0.5592"
$ws.Range("E47").Value = "  -4.86%  "
$ws.Range("D48").Value = "'1.904"
$ws.Range("E48").Value = "  -3.22%  "
$ws.Range("D49").Value = "'118.28"
$ws.Range("E49").Value = "  -3.33%  "
$ws.Range("D50").Value = "'1.136"
$ws.Range("E50").Value = "  -3.67%  "
$ws.Range("D51").Value = "'0.06820"
$ws.Range("E51").Value = "  -3.36%  "
